$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary" -- update evaluation metrics for one_svm/split_5/test_50_50
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.649812734082397
$wsSummary.Range("C2").Value = 0.590702947845805
$wsSummary.Range("D2").Value = 0.9756554307116105
$wsSummary.Range("E2").Value = 0.7358757062146892
$wsSummary.Range("F2").Value = 0.8631544068919814
$wsSummary.Range("G2").Value = 0.9517987633501968
$wsSummary.Range("H2").Value = 0.7969216849724361
$wsSummary.Range("I2").Value = 521
$wsSummary.Range("J2").Value = 361
$wsSummary.Range("K2").Value = 173
$wsSummary.Range("L2").Value = 13

# ---------------------------------------------------------------------------
# Sheet "Classification Report"
# ---------------------------------------------------------------------------
$wsReport = $wb.Worksheets.Item("Classification Report")

# row 2 -> class "0"
$wsReport.Range("B2").Value = 0.9301075268817204
$wsReport.Range("C2").Value = 0.3239700374531835
$wsReport.Range("D2").Value = 0.4805555555555556

# row 3 -> class "1"
$wsReport.Range("B3").Value = 0.590702947845805
$wsReport.Range("C3").Value = 0.9756554307116105
$wsReport.Range("D3").Value = 0.7358757062146892

# row 4 -> accuracy
$wsReport.Range("B4").Value = 0.649812734082397
$wsReport.Range("C4").Value = 0.649812734082397
$wsReport.Range("D4").Value = 0.649812734082397
$wsReport.Range("E4").Value = 0.649812734082397

# row 5 -> macro avg
$wsReport.Range("B5").Value = 0.7604052373637626
$wsReport.Range("C5").Value = 0.649812734082397
$wsReport.Range("D5").Value = 0.6082156308851224

# row 6 -> weighted avg
$wsReport.Range("B6").Value = 0.7604052373637626
$wsReport.Range("C6").Value = 0.649812734082397
$wsReport.Range("D6").Value = 0.6082156308851224

# ---------------------------------------------------------------------------
# Sheet "Confusion Matrix"
# ---------------------------------------------------------------------------
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")

# row 2 -> Actual 0
$wsConfusion.Range("B2").Value = 173
$wsConfusion.Range("C2").Value = 361

# row 3 -> Actual 1
$wsConfusion.Range("B3").Value = 13
$wsConfusion.Range("C3").Value = 521
